$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data row 32 (Pipe Diameters table extension) -----------------
# Copy the formatting (style) from the last existing row of the table so the
# new row visually matches the rest of the table before we populate values.
$ws.Range("B31").Copy()
$ws.Range("B32").PasteSpecial(-4122)  # xlPasteFormats

# New input value for A32 (kept without special formatting, as in source)
$ws.Range("A32").Value = 2600

# Re-apply the shared formula across the whole B21:B32 range so that the
# existing rows keep their computed values and the newly added B32 cell
# picks up the same relative formula, extended one row further down.
$ws.Range("B21:B32").Formula = '=(A21 * 1000000) / ($E$20 * 1000 * $E$21)'

# --- Update the visible selection to the newly added cell -----------------
$ws.Activate()
$ws.Range("A32").Select()
